$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (column G) values computed from the regenerated std/mean-based s_vals calc.
$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 2
    6  = 2
    7  = 1
    8  = 1
    9  = 3
    10 = 0
    11 = 0
    12 = 0
    14 = 1
    15 = 3
    16 = 2
    17 = 3
    18 = 2
    19 = 1
    20 = 0
    21 = 1
    22 = 2
    23 = 1
    24 = 2
    25 = 2
    26 = 1
    27 = 1
    28 = 2
    30 = 1
    32 = 2
    33 = 2
    34 = 2
    35 = 1
    36 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
